$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the expiry date value in D7 (45787 -> 45797)
$ws.Range("D7").Value = 45797

# Update the selected/active cell shown in the sheet view (D4 -> E7)
$ws.Range("E7").Select()
